$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed price values in column D that are plain numbers (e.g. "246.09")
# must stay text (to preserve exact formatting such as trailing zeros), so those
# cells are switched to a Text number format before the new value is written.
# (Cells whose text contains two "." separators, e.g. "42.405.42", are never
# auto-converted to numbers, so they do not need this treatment.)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = '42.405.42'
$ws.Range("E2").Value = '  +0.68%  '
$ws.Range("D3").Value = '2.247.58'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '246.09'
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("D6").Value = '0.631'
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("D7").Value = '75.93'
$ws.Range("E7").Value = '  -0.96%  '
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("D9").Value = '0.618'
$ws.Range("E9").Value = '  -1.48%  '
$ws.Range("D10").Value = '44.40'
$ws.Range("E10").Value = '  +7.44%  '
$ws.Range("D11").Value = '0.0948'
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("D12").Value = '7.23'
$ws.Range("E12").Value = '  +0.95%  '
$ws.Range("D13").Value = '0.103'
$ws.Range("E13").Value = '  -1.12%  '
$ws.Range("D14").Value = '14.58'
$ws.Range("E14").Value = '  -1.88%  '
$ws.Range("D15").Value = '0.858'
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").Value = '2.245.76'
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("D17").Value = '42.274.30'
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("E18").Value = '  +3.79%  '
$ws.Range("D19").Value = '6.19'
$ws.Range("E19").Value = '  +0.60%  '
$ws.Range("D20").Value = '72.21'
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("E21").Value = '  +4.23%  '
$ws.Range("D22").Value = '231.79'
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("D23").Value = '8.92'
$ws.Range("E23").Value = '  +22.06%  '
$ws.Range("D25").Value = '11.48'
$ws.Range("E25").Value = '  +2.65%  '
$ws.Range("D26").Value = '3.63'
$ws.Range("E26").Value = '  -2.42%  '
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("E28").Value = '  +1.73%  '
$ws.Range("E29").Value = '  -1.13%  '
$ws.Range("D30").Value = '20.68'
$ws.Range("E30").Value = '  +0.85%  '
$ws.Range("E31").Value = '  -3.40%  '
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").Value = '30.80'
$ws.Range("E33").Value = '  -5.73%  '
$ws.Range("D34").Value = '5.36'
$ws.Range("E34").Value = '  +10.34%  '
$ws.Range("D35").Value = '0.125'
$ws.Range("E35").Value = '  -0.38%  '
$ws.Range("D36").Value = '4.57'
$ws.Range("E36").Value = '  -0.87%  '
$ws.Range("D37").Value = '0.0315'
$ws.Range("E37").Value = '  +6.14%  '
$ws.Range("D38").Value = '13.99'
$ws.Range("E38").Value = '  +6.25%  '
$ws.Range("D39").Value = '2.18'
$ws.Range("E39").Value = '  -1.13%  '
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("D41").Value = '63.80'
$ws.Range("E41").Value = '  +6.05%  '
$ws.Range("D43").Value = '107.64'
$ws.Range("E43").Value = '  -6.05%  '
$ws.Range("D44").Value = '8.79'
$ws.Range("E44").Value = '  -0.33%  '
$ws.Range("E45").Value = '  +2.52%  '
$ws.Range("D46").Value = '0.998'
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("E47").Value = '  +0.39%  '
$ws.Range("E48").Value = '  +2.08%  '
$ws.Range("D49").Value = '2.37'
$ws.Range("E49").Value = '  +4.23%  '
$ws.Range("E50").Value = '  +0.77%  '
$ws.Range("E51").Value = '  +0.97%  '
